$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '41.287.68'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -1.69%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.184.31'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -1.15%  '

$ws.Range("E4").Value = '  +0.00%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '238.07'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.87%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.614'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -1.62%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '70.23'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -3.66%  '

$ws.Range("E8").Value = '  +0.00%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.576'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -4.31%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '40.10'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -4.84%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0927'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -2.40%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '54.57'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -5.00%  '

$ws.Range("E13").Value = '  -1.59%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.77'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -3.95%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '2.508.30'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -1.38%  '

$ws.Range("E16").Value = '  -0.78%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '2.182.54'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -1.28%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.803'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -3.91%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '41.176.20'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -1.66%  '

$ws.Range("E20").Value = '  -5.52%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '70.83'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -2.20%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.94'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -2.56%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '9.73'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -3.60%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '226.90'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -0.82%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '1.93'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -6.11%  '

$ws.Range("E26").Value = '  +0.07%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '10.84'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -6.67%  '

$ws.Range("E28").Value = '  -1.77%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.21'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -2.01%  '

$ws.Range("E30").Value = '  +0.40%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '167.61'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +0.40%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '20.04'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -2.26%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '30.91'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +7.79%  '

$ws.Range("E34").Value = '  -2.10%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '5.15'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -7.87%  '

$ws.Range("E36").Value = '  -2.65%  '

$ws.Range("E37").Value = '  -6.33%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '4.12'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -2.71%  '

$ws.Range("E39").Value = '  -4.63%  '

$ws.Range("E40").Value = '  -1.15%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '11.73'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -11.18%  '

$ws.Range("E42").Value = '  -2.78%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '59.65'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -7.49%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.191'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -2.86%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0980'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -2.05%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '8.28'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -4.62%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '98.58'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -4.76%  '

$ws.Range("E48").Value = '  -1.24%  '

$ws.Range("E49").Value = '  -1.95%  '

$ws.Range("E50").Value = '  -6.26%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.62'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -2.51%  '
